# Applies the "resolve domain related templates" edit described by the diff.

$wb = $excel.ActiveWorkbook

# --- Sheet "Menus": add two new rows (MenuItem1 / DataDict1 domains) ---
$wsMenus = $wb.Worksheets.Item("Menus")
$wsMenus.Range("A5").Value = "MenuItem1"
$wsMenus.Range("B5").Value = "/menu/MenuItem1"
$wsMenus.Range("C5").Value = "ADMINISTRATION"
$wsMenus.Range("D5").Value = "MAINTENANCE"

$wsMenus.Range("A6").Value = "DataDict1"
$wsMenus.Range("B6").Value = "/dict/DataDict1"
$wsMenus.Range("C6").Value = "ADMINISTRATION"
$wsMenus.Range("D6").Value = "MAINTENANCE"

$wsMenus.Range("D8").Select()

# --- Sheet "Domains": rename the MenuItem / DataDict class fields ---
$wsDomains = $wb.Worksheets.Item("Domains")
$wsDomains.Range("B3").Value = "MenuItem1"
$wsDomains.Range("B8").Value = "DataDict1"
$wsDomains.Range("B13").Select()

# --- Sheet "System Properties": fix spring.mvc.locale value ---
$wsSysProps = $wb.Worksheets.Item("System Properties")
$wsSysProps.Range("B7").Value = "zh_CN"
$wsSysProps.Range("B8").Select()

# --- Sheet "Messages": scroll the view (no content change) ---
$wsMessages = $wb.Worksheets.Item("Messages")
$wsMessages.Range("B20").Select()

# --- Activate "Menus" tab (was "Domains") ---
$wsMenus.Activate()
$wsMenus.Select()
